# Applies the author's edits described by the commit "minimax and combat modeling":
#   - Rename sheet "Comunication" -> "Movement+Comunication"
#   - Update each sheet's last-used selection / active cell
#   - Leave "Class26AUG" as the final active sheet/tab

$wb = $excel.ActiveWorkbook

# --- Rename the "Comunication" sheet -------------------------------------
$wsComunication = $wb.Sheets.Item("Comunication")
$wsComunication.Name = "Movement+Comunication"

# --- Update per-sheet selections ------------------------------------------

# Movement+Comunication: selection I13 -> J12
$wsComunication.Activate()
$wsComunication.Range("J12").Select()

# HW7: selection D24 -> H5
$wsHW7 = $wb.Sheets.Item("HW7")
$wsHW7.Activate()
$wsHW7.Range("H5").Select()

# Class24AUG: selection J18 -> B11
$wsClass24 = $wb.Sheets.Item("Class24AUG")
$wsClass24.Activate()
$wsClass24.Range("B11").Select()

# Sheet2: selection stays E21 (re-select so state is explicit)
$wsSheet2 = $wb.Sheets.Item("Sheet2")
$wsSheet2.Activate()
$wsSheet2.Range("E21").Select()

# Class26AUG: selection X7 -> U2; ends up the active tab
$wsClass26 = $wb.Sheets.Item("Class26AUG")
$wsClass26.Activate()

# Re-assert the (already-centered) alignment on T1:U1 so the style table
# collapses the redundant "centered" cell format onto the one also used
# elsewhere on the sheet (matches the cellXfs cleanup in the saved file).
$wsClass26.Range("T1:U1").HorizontalAlignment = -4108

$wsClass26.Range("U2").Select()
